$d = $word.ActiveDocument

function Find-ParagraphByText($doc, $exactText) {
    for ($i = 1; $i -le $doc.Paragraphs.Count; $i++) {
        $p = $doc.Paragraphs.Item($i)
        if ($p.Range.Text -eq $exactText) {
            return $i
        }
    }
    return -1
}

# ---------------------------------------------------------------
# 1. Education date: "... (Spring 2015)" -> "... (Spring 2016)"
# ---------------------------------------------------------------
$r = $d.Content
$r.Find.Execute("2015)", $true, $false, $false, $false, $false, $true, 1, $false, "2016)", 2) | Out-Null

# ---------------------------------------------------------------
# 2. Job title: "Software Development Intern" -> "Web Applications Developer"
# ---------------------------------------------------------------
$r = $d.Content
$r.Find.Execute("Software Development Intern", $true, $false, $false, $false, $false, $true, 1, $false, "Web Applications Developer", 2) | Out-Null

# ---------------------------------------------------------------
# 3. Skills bullets: merge "Middleware: AJAX, JSON" + "Back end: PHP, Python"
#    into a single "Back end: PHP, Python, Tornado" bullet, removing the
#    old separate "Back end: PHP, Python" bullet paragraph.
# ---------------------------------------------------------------
$idx = Find-ParagraphByText $d "Middleware: AJAX, JSON`r"
if ($idx -ne -1) {
    $p = $d.Paragraphs.Item($idx)
    $rr = $p.Range
    $rr.MoveEnd(1, -1) | Out-Null
    $rr.Text = "Back end: PHP, Python, Tornado"
}

$idx2 = Find-ParagraphByText $d "Back end: PHP, Python`r"
if ($idx2 -ne -1) {
    $p2 = $d.Paragraphs.Item($idx2)
    $p2.Range.Delete() | Out-Null
}

# ---------------------------------------------------------------
# 4. StarChase bullet: "Implement changes to company website" ->
#    "Maintain company website (http://starchase.com/)" and add a new
#    bullet after it for the Pursuit for Change Wordpress site.
# ---------------------------------------------------------------
$r = $d.Content
$r.Find.Execute("Implement changes to company website", $true, $false, $false, $false, $false, $true, 1, $false, "Maintain company website (http://starchase.com/)", 2) | Out-Null

$idx3 = Find-ParagraphByText $d "Maintain company website (http://starchase.com/)`r"
if ($idx3 -ne -1) {
    $p3 = $d.Paragraphs.Item($idx3)
    $p3.Range.InsertParagraphAfter() | Out-Null
    $newPara = $d.Paragraphs.Item($idx3 + 1)
    $newPara.Range.Text = "Develop wordpress website for Pursuit for Change (https://www.pursuitforchange.org/)"
}

# ---------------------------------------------------------------
# 5. Skills: "Languages: JavaScript, PHP, Python, Java" -> drop ", Java"
# ---------------------------------------------------------------
$r = $d.Content
$r.Find.Execute(" Python, Java", $true, $false, $false, $false, $false, $true, 1, $false, " Python", 2) | Out-Null

# ---------------------------------------------------------------
# 6. Tools: "Sublime Text 2" -> "Sublime Text 3"
# ---------------------------------------------------------------
$r = $d.Content
$r.Find.Execute("Sublime Text 2", $true, $false, $false, $false, $false, $true, 1, $false, "Sublime Text 3", 2) | Out-Null

Write-Host "Done."
